# Refresh the cryptos price/volume snapshot (column D = Price, column E =
# Volume(1h)) for the rows whose figures moved since the last run.
#
# Values that read as plain numbers (e.g. "1.00", "555.49") are written with
# a leading apostrophe so Excel stores them as literal text -- matching the
# rest of the sheet, where prices/percentages are kept as text (e.g.
# "3.380.77", "  +0.34%  ") rather than numeric values, and so trailing
# zeros (e.g. "18.30" vs "18.3") are preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "65.097.14"
$ws.Cells.Item(2, 5).Value = "  +0.76%  "
$ws.Cells.Item(3, 4).Value = "3.380.77"
$ws.Cells.Item(3, 5).Value = "  +0.34%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).Value = "'555.49"
$ws.Cells.Item(5, 5).Value = "  -0.35%  "
$ws.Cells.Item(6, 4).Value = "'174.75"
$ws.Cells.Item(6, 5).Value = "  -0.87%  "
$ws.Cells.Item(7, 4).Value = "'0.633"
$ws.Cells.Item(7, 5).Value = "  +1.98%  "
$ws.Cells.Item(8, 4).Value = "3.370.97"
$ws.Cells.Item(8, 5).Value = "  +0.39%  "
$ws.Cells.Item(9, 5).Value = "  -0.07%  "
$ws.Cells.Item(10, 5).Value = "  +5.86%  "
$ws.Cells.Item(11, 5).Value = "  +1.09%  "
$ws.Cells.Item(12, 4).Value = "'53.62"
$ws.Cells.Item(12, 5).Value = "  -1.84%  "
$ws.Cells.Item(13, 4).Value = "'0.0000279"
$ws.Cells.Item(13, 5).Value = "  +1.95%  "
$ws.Cells.Item(14, 4).Value = "'9.18"
$ws.Cells.Item(14, 5).Value = "  +0.96%  "
$ws.Cells.Item(15, 4).Value = "3.920.38"
$ws.Cells.Item(15, 5).Value = "  +0.43%  "
$ws.Cells.Item(16, 4).Value = "'18.30"
$ws.Cells.Item(17, 4).Value = "3.390.93"
$ws.Cells.Item(17, 5).Value = "  +0.61%  "
$ws.Cells.Item(18, 5).Value = "  -0.40%  "
$ws.Cells.Item(19, 4).Value = "'11.90"
$ws.Cells.Item(19, 5).Value = "  +0.13%  "
$ws.Cells.Item(20, 4).Value = "64.947.15"
$ws.Cells.Item(20, 5).Value = "  +0.68%  "
$ws.Cells.Item(21, 4).Value = "'0.999"
$ws.Cells.Item(21, 5).Value = "  +1.27%  "
$ws.Cells.Item(22, 5).Value = "  -1.12%  "
$ws.Cells.Item(23, 5).Value = "  +0.45%  "
$ws.Cells.Item(24, 4).Value = "'14.36"
$ws.Cells.Item(24, 5).Value = "  +7.51%  "
$ws.Cells.Item(25, 5).Value = "  -0.59%  "
$ws.Cells.Item(26, 4).Value = "'87.51"
$ws.Cells.Item(26, 5).Value = "  +2.44%  "
$ws.Cells.Item(27, 5).Value = "  +0.96%  "
$ws.Cells.Item(28, 4).Value = "'10.71"
$ws.Cells.Item(28, 5).Value = "  -2.29%  "
$ws.Cells.Item(29, 4).Value = "'8.71"
$ws.Cells.Item(29, 5).Value = "  -1.25%  "
$ws.Cells.Item(30, 4).Value = "'31.19"
$ws.Cells.Item(30, 5).Value = "  +3.91%  "
$ws.Cells.Item(31, 4).Value = "'6.53"
$ws.Cells.Item(31, 5).Value = "  -1.24%  "
$ws.Cells.Item(32, 4).Value = "'63.01"
$ws.Cells.Item(32, 5).Value = "  +7.21%  "
$ws.Cells.Item(33, 4).Value = "'11.47"
$ws.Cells.Item(33, 5).Value = "  -0.38%  "
$ws.Cells.Item(34, 4).Value = "'577.90"
$ws.Cells.Item(34, 5).Value = "  -1.19%  "
$ws.Cells.Item(35, 5).Value = "  -0.42%  "
$ws.Cells.Item(37, 4).Value = "'3.65"
$ws.Cells.Item(37, 5).Value = "  +4.85%  "
$ws.Cells.Item(38, 5).Value = "  +0.83%  "
$ws.Cells.Item(39, 4).Value = "'35.70"
$ws.Cells.Item(39, 5).Value = "  -0.42%  "
$ws.Cells.Item(40, 5).Value = "  -0.23%  "
$ws.Cells.Item(41, 4).Value = "0.0₃0739"
$ws.Cells.Item(41, 5).Value = "  -2.60%  "
$ws.Cells.Item(42, 4).Value = "3.099.80"
$ws.Cells.Item(42, 5).Value = "  -0.39%  "
$ws.Cells.Item(43, 4).Value = "'0.0417"
$ws.Cells.Item(43, 5).Value = "  +1.43%  "
$ws.Cells.Item(44, 4).Value = "'2.77"
$ws.Cells.Item(44, 5).Value = "  -1.59%  "
$ws.Cells.Item(45, 4).Value = "'2.45"
$ws.Cells.Item(45, 5).Value = "  -2.04%  "
$ws.Cells.Item(46, 5).Value = "  +2.47%  "
$ws.Cells.Item(47, 4).Value = "'3.17"
$ws.Cells.Item(47, 5).Value = "  -1.67%  "
$ws.Cells.Item(48, 4).Value = "'1.00"
$ws.Cells.Item(48, 5).Value = "  +0.10%  "
$ws.Cells.Item(49, 4).Value = "'141.84"
$ws.Cells.Item(49, 5).Value = "  +5.04%  "
$ws.Cells.Item(50, 5).Value = "  -2.33%  "
$ws.Cells.Item(51, 4).Value = "'8.29"
$ws.Cells.Item(51, 5).Value = "  -0.57%  "
